$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: "On progress" -> "Done"
$ws.Range("B7").Value = "Done"

# New row 8: Exercise = "Scene" (reuses/renames the now-unused "On progress" shared string)
$ws.Range("A8").Value = "Scene"

# Expand the table to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B8"))

# Expand the color-scale conditional formatting range to include the new row
$fc = $ws.Range("B2:B7").FormatConditions.Item(2)
$fc.ModifyAppliesToRange($ws.Range("B2:B8"))

# Move the active selection to the new last cell
$ws.Range("B8").Select()
